$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '63.622.82'
$ws.Range('E2').Value = '  -1.08%  '
$ws.Range('D3').Value = '3.406.39'
$ws.Range('E3').Value = '  -0.16%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '568.39'
$ws.Range('E5').Value = '  -0.26%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '156.28'
$ws.Range('E7').Value = '  +0.01%  '
$ws.Range('D8').Value = '3.408.15'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.570'
$ws.Range('E9').Value = '  -6.73%  '
$ws.Range('E10').Value = '  +0.78%  '
$ws.Range('E11').Value = '  -2.40%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.426'
$ws.Range('E12').Value = '  -3.17%  '
$ws.Range('D13').Value = '3.999.78'
$ws.Range('E13').Value = '  +0.00%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '27.03'
$ws.Range('E15').Value = '  -2.77%  '
$ws.Range('E16').Value = '  -7.28%  '
$ws.Range('D17').Value = '63.732.94'
$ws.Range('E17').Value = '  -1.07%  '
$ws.Range('D18').Value = '3.406.53'
$ws.Range('E18').Value = '  +0.37%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.08'
$ws.Range('E19').Value = '  -3.89%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '13.55'
$ws.Range('E20').Value = '  -2.72%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '381.09'
$ws.Range('E21').Value = '  +1.93%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '7.74'
$ws.Range('E22').Value = '  -2.63%  '
$ws.Range('E23').Value = '  -0.03%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '71.01'
$ws.Range('E24').Value = '  -1.61%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.517'
$ws.Range('E25').Value = '  -5.73%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.0000114'
$ws.Range('E26').Value = '  -3.73%  '
$ws.Range('E27').Value = '  -3.79%  '
$ws.Range('E28').Value = '  +0.93%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.996'
$ws.Range('E29').Value = '  -0.36%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '6.06'
$ws.Range('E30').Value = '  -0.65%  '
$ws.Range('E31').Value = '  -6.24%  '
$ws.Range('E32').Value = '  -1.24%  '
$ws.Range('B33').Value = 'USDe'
$ws.Range('C33').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.998'
$ws.Range('E33').Value = '  -0.05%  '
$ws.Range('B34').Value = 'EthereumClassic'
$ws.Range('C34').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '22.91'
$ws.Range('E34').Value = '  -0.62%  '
$ws.Range('B35').Value = 'Aptos'
$ws.Range('C35').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '6.96'
$ws.Range('E35').Value = '  -3.33%  '
$ws.Range('B36').Value = 'ImmutableX'
$ws.Range('C36').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.52'
$ws.Range('E36').Value = '  -4.87%  '
$ws.Range('B37').Value = 'Monero'
$ws.Range('C37').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '161.05'
$ws.Range('E37').Value = '  +0.35%  '
$ws.Range('B38').Value = 'Mantle'
$ws.Range('C38').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.831'
$ws.Range('E38').Value = '  +8.42%  '
$ws.Range('B39').Value = 'Stacks'
$ws.Range('C39').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.83'
$ws.Range('E39').Value = '  -2.50%  '
$ws.Range('B40').Value = 'EnergySwap'
$ws.Range('C40').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '26.09'
$ws.Range('E40').Value = '  -2.15%  '
$ws.Range('B41').Value = 'Maker'
$ws.Range('C41').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D41').Value = '2.802.51'
$ws.Range('E41').Value = '  -1.57%  '
$ws.Range('B42').Value = 'Hedera'
$ws.Range('C42').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0722'
$ws.Range('E42').Value = '  -4.53%  '
$ws.Range('B43').Value = 'OKB'
$ws.Range('C43').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '42.77'
$ws.Range('E43').Value = '  +0.20%  '
$ws.Range('B44').Value = 'RenderToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '6.37'
$ws.Range('E44').Value = '  -6.08%  '
$ws.Range('B45').Value = 'Filecoin'
$ws.Range('C45').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '4.39'
$ws.Range('E45').Value = '  -4.66%  '
$ws.Range('B46').Value = 'InjectiveProtocol'
$ws.Range('C46').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '25.57'
$ws.Range('E46').Value = '  -1.89%  '
$ws.Range('B47').Value = 'VeChain'
$ws.Range('C47').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0304'
$ws.Range('E47').Value = '  -2.61%  '
$ws.Range('B48').Value = 'Bittensor'
$ws.Range('C48').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '325.35'
$ws.Range('E48').Value = '  +3.07%  '
$ws.Range('B49').Value = 'dogwifhat'
$ws.Range('C49').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.30'
$ws.Range('E49').Value = '  +7.82%  '
$ws.Range('B50').Value = 'ONDO'
$ws.Range('C50').Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.03'
$ws.Range('E50').Value = '  -3.80%  '
$ws.Range('B51').Value = 'Stellar'
$ws.Range('C51').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.103'
$ws.Range('E51').Value = '  -5.13%  '
